$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so Excel does not
# auto-convert numeric-looking / percentage-looking strings into
# real numbers, dates, etc. The source workbook stores every one
# of these as an inline/shared string (t="inlineStr"), never a
# true number, so we must preserve that by using NumberFormat "@".

$ws.Range("D2,E2,D3,E3,E4,D5,E5,D6,E6,E7,D8,E8,D9,E9,E10,D11,E11,D12,E12,E13,D14,E14,D15,E15,D16,E16,D17,E17,D18,E18,D19,E19,D20,E20,D21,E21,D22,E22,D23,E23,D24,E24,D25,E25,E26,E27,E28,D29,E29,D30,E30,E31,D32,E32,D33,E33,D34,E34,D35,E35,E36,B37,C37,D37,E37,B38,C38,D38,E38,B39,C39,D39,E39,B40,C40,D40,E40,B41,C41,D41,E41,B42,C42,D42,E42,B43,C43,D43,E43,D44,E44,B45,C45,D45,E45,B46,C46,D46,E46,B47,C47,D47,E47,B48,C48,D48,E48,B49,C49,D49,E49,B50,C50,D50,E50,B51,C51,D51,E51").NumberFormat = "@"

$ws.Range('D2').Value = '55.477.66'
$ws.Range('E2').Value = '  -5.85%  '
$ws.Range('D3').Value = '2.925.94'
$ws.Range('E3').Value = '  -9.43%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '467.06'
$ws.Range('E5').Value = '  -13.03%  '
$ws.Range('D6').Value = '123.42'
$ws.Range('E6').Value = '  -9.23%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '2.920.12'
$ws.Range('E8').Value = '  -9.60%  '
$ws.Range('D9').Value = '0.398'
$ws.Range('E9').Value = '  -13.23%  '
$ws.Range('E10').Value = '  -13.27%  '
$ws.Range('D11').Value = '0.0947'
$ws.Range('E11').Value = '  -17.87%  '
$ws.Range('D12').Value = '0.327'
$ws.Range('E12').Value = '  -17.27%  '
$ws.Range('E13').Value = '  -3.38%  '
$ws.Range('D14').Value = '3.434.33'
$ws.Range('E14').Value = '  -9.18%  '
$ws.Range('D15').Value = '22.45'
$ws.Range('E15').Value = '  -13.86%  '
$ws.Range('D16').Value = '55.704.49'
$ws.Range('E16').Value = '  -5.66%  '
$ws.Range('D17').Value = '2.925.64'
$ws.Range('E17').Value = '  -9.45%  '
$ws.Range('D18').Value = '0.0000132'
$ws.Range('E18').Value = '  -17.34%  '
$ws.Range('D19').Value = '5.09'
$ws.Range('E19').Value = '  -13.78%  '
$ws.Range('D20').Value = '11.45'
$ws.Range('E20').Value = '  -13.42%  '
$ws.Range('D21').Value = '7.01'
$ws.Range('E21').Value = '  -15.36%  '
$ws.Range('D22').Value = '306.37'
$ws.Range('E22').Value = '  -15.20%  '
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').Value = '0.445'
$ws.Range('E24').Value = '  -14.40%  '
$ws.Range('D25').Value = '58.85'
$ws.Range('E25').Value = '  -16.78%  '
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('E27').Value = '  -9.56%  '
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').Value = '0.0₃0800'
$ws.Range('E29').Value = '  -17.74%  '
$ws.Range('D30').Value = '5.87'
$ws.Range('E30').Value = '  -16.84%  '
$ws.Range('E31').Value = '  -9.68%  '
$ws.Range('D32').Value = '19.05'
$ws.Range('E32').Value = '  -13.58%  '
$ws.Range('D33').Value = '6.02'
$ws.Range('E33').Value = '  -15.03%  '
$ws.Range('D34').Value = '1.57'
$ws.Range('E34').Value = '  -18.54%  '
$ws.Range('D35').Value = '144.50'
$ws.Range('E35').Value = '  -11.78%  '
$ws.Range('E36').Value = '  -15.95%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = '5.35'
$ws.Range('E37').Value = '  -16.00%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '1.21'
$ws.Range('E38').Value = '  -15.86%  '
$ws.Range('B39').Value = 'RenzoRestakedETH'
$ws.Range('C39').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D39').Value = '2.959.40'
$ws.Range('E39').Value = '  -9.26%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').Value = '0.0605'
$ws.Range('E41').Value = '  -14.56%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '21.40'
$ws.Range('E42').Value = '  -17.29%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = '34.86'
$ws.Range('E43').Value = '  -15.11%  '
$ws.Range('D44').Value = '0.954'
$ws.Range('E44').Value = '  -12.58%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '0.598'
$ws.Range('E45').Value = '  -16.64%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '3.38'
$ws.Range('E46').Value = '  -15.68%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.052.18'
$ws.Range('E47').Value = '  -10.37%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = '1.29'
$ws.Range('E48').Value = '  -14.07%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = '5.25'
$ws.Range('E49').Value = '  -16.30%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '17.49'
$ws.Range('E50').Value = '  -15.83%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = '0.0208'
$ws.Range('E51').Value = '  -13.97%  '
